$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.455362044514542, 10.34677158129881, 261.3203778131603, 10.19245300693656, 1, 283.3149644459102)
    3  = @(0.2917716402565462, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 2.591208233317391)
    4  = @(0.04271373187048222, 1.655778082260271, 0.1494219747398047, 10.19245300693656, 1, 12.04036679580711)
    5  = @(0.1190320826869504, 117.745847958593, 0.7527432677738641, 10.19245300693656, 1, 128.8100763159904)
    6  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 4.358119930609447)
    7  = @(0.6606524410359556, 0.04071648406533734, 22.3905356188092, 10.19245300693656, 0, 33.28435755084705)
    8  = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 7.143138311642302)
    9  = @(0.1190320826869504, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1, 1.406728370586922)
    10 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
